$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SGS")

# Insert a new row at the top for the "Grit-O" section header
$ws.Rows.Item(1).Insert()
$ws.Range("B1").Value = "Grit-O"
$ws.Range("B1").HorizontalAlignment = -4108
$ws.Range("C1").HorizontalAlignment = -4108
$ws.Range("D1").HorizontalAlignment = -4108
$ws.Range("B1:D1").Merge()

# Rows 2-4 now hold the original EO-style header + 2 data rows (unchanged values)
# Append the new "Grit-S" section starting at row 5
$ws.Range("B5").Value = "Grit-S"
$ws.Range("B5").HorizontalAlignment = -4108
$ws.Range("C5").HorizontalAlignment = -4108
$ws.Range("D5").HorizontalAlignment = -4108
$ws.Range("B5:D5").Merge()

$ws.Range("B6").Value = "Sample Units"
$ws.Range("C6").Value = "items"
$ws.Range("D6").Value = "A-C"

$ws.Range("A7").Value = "Consistency of interest"
$ws.Range("B7").Value = 188
$ws.Range("C7").Value = 4
$ws.Range("D7").Value = 0.674

$ws.Range("A8").Value = "perseverance of efforts"
$ws.Range("B8").Value = 188
$ws.Range("C8").Value = 4
$ws.Range("D8").Value = 0.604

Write-Host "done"
